$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for API ID / API Secret (row 2 already uses customFormat style 1)
$ws.Range("T2").Value = "API ID"
$ws.Range("U2").Value = "API Secret"

# Create a new row 3 with a cell (F3) styled like a hyperlink (adds Hyperlink cell style/font),
# then remove the actual hyperlink object so only the formatting remains.
$ws.Hyperlinks.Add($ws.Range("F3"), "", "", "", "") | Out-Null
$ws.Hyperlinks.Delete()

# Update the selection/active cell shown in the view (also clears the old topLeftCell scroll state)
$ws.Range("A3:Y3").Select()
